$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Generator Data")
$ws1.Range("B2").Value = 32708.748857899998
$ws1.Range("C2").Value = 27546.739683299998
$ws1.Range("B3").Value = 13740.945395203789
$ws1.Range("C3").Value = 8264.0219049899988
$ws1.Range("B4").Value = 1374.0945395203789
$ws1.Range("C4").Value = 743.76197144909986
$ws1.Range("B5").Value = 122896.074876
$ws1.Range("C5").Value = 45825.056022199999

$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")
$ws2.Range("B2").Value = 84311.981381491467
$ws2.Range("B3").Value = 84409.149516779638
